$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C, rows 2 through 339 all hold the date serial 45189 ("Förändrad" column).
# Update them to 45190 (one day later) to match the target edit.
$ws.Range("C2:C339").Value = 45190
